$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1466
$ws1.Range("F3").Value = 3082
$ws1.Range("F4").Value = 42
$ws1.Range("F5").Value = 670
$ws1.Range("F6").Value = 290

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1466
$ws4.Range("F3").Value = 3082
$ws4.Range("F4").Value = 42
$ws4.Range("F5").Value = 670
$ws4.Range("F7").Value = 290
